$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Updated Neo4j query text for the three tabs: each query gets an
# "order by ... LIMIT 100" clause appended (per commit "Fixed Bento 80 Test scripts"). ---

$newCasesTab = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss:study_subject)
	WHERE ss.study_subject_id in ['BENTO-CASE-3282798','BENTO-CASE-3292831','BENTO-CASE-3295670','BENTO-CASE-3295756','BENTO-CASE-3296613','BENTO-CASE-4187183','BENTO-CASE-4213980']
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
RETURN DISTINCT 
	ss.study_subject_id AS `Case ID`,
	p.program_acronym AS `Program Code`,
	p.program_id AS `Program ID`,
	s.study_acronym AS `Arm`,
	ss.disease_subtype AS `Diagnosis`,
	sf.grouped_recurrence_score AS `Recurrence Score`,
	d.tumor_size_group AS `Tumor Size (cm)`,
	d.er_status AS `ER Status`,
	d.pr_status AS `PR Status`,
	demo.age_at_index AS `Age (years)`,
	demo.survival_time AS `Survival (days)`
 order By ss.study_subject_id ASC LIMIT 100
'@

$newSamplesTab = @'
MATCH (ss:study_subject)
	WHERE ss.study_subject_id in ['BENTO-CASE-3282798','BENTO-CASE-3292831','BENTO-CASE-3295670','BENTO-CASE-3295756','BENTO-CASE-3296613','BENTO-CASE-4187183','BENTO-CASE-4213980']
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[*..2]-(parent)<--(f:file)
OPTIONAL MATCH (f)-[:file_of_laboratory_procedure]->(lp)
RETURN DISTINCT 
	samp.sample_id AS `Sample ID`,
	ss.study_subject_id AS `Case ID`,
	p.program_acronym AS `Program Code`,
	s.study_acronym AS `Arm`,
	ss.disease_subtype AS `Diagnosis`,
	samp.tissue_type AS `Tissue Type`,
	samp.composition AS `Tissue Composition`,
	samp.sample_anatomic_site AS `Sample Anatomic Site`,
	samp.method_of_sample_procurement AS `Sample Procurement Method`,
	lp.test_name AS `platform`
 order By samp.sample_id ASC LIMIT 100
'@

$newFilesTab = @'
MATCH (ss:study_subject)
	WHERE ss.study_subject_id in ['BENTO-CASE-3282798','BENTO-CASE-3292831','BENTO-CASE-3295670','BENTO-CASE-3295756','BENTO-CASE-3296613','BENTO-CASE-4187183','BENTO-CASE-4213980']
MATCH (ss)<-[*..2]-(parent)<--(f:file)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
RETURN DISTINCT 
	f.file_name AS `File Name`,
	head(labels(parent)) AS `Association`,
	f.file_description AS `Description`,
	f.file_format AS `File Format`,
	f.file_size AS `Size`,
	p.program_acronym AS `Program Code`,
	s.study_acronym AS `Arm`,
	ss.study_subject_id AS `Case ID`,
	samp.sample_id AS `Sample ID`
 order By f.file_name ASC LIMIT 100
'@

$ws.Range("B2").Value = $newCasesTab
$ws.Range("B3").Value = $newSamplesTab
$ws.Range("B4").Value = $newFilesTab

# --- Row heights grow slightly because each query now wraps one extra line. ---
$ws.Rows.Item(2).RowHeight = 374.4
$ws.Rows.Item(3).RowHeight = 288
$ws.Rows.Item(4).RowHeight = 273.6

# --- Selection moved from C4 back up to C3, scrolled back to the top of the sheet. ---
[void]$ws.Activate()
[void]$ws.Range("C3").Select()
